# ---------------------------------------------------------------------------
# Applies the "Add files via upload" revision to the Code of Conduct document:
#   1. "Group J" (team-name placeholder) -> "C#"
#   2. Fill in the blank "Date" cell of the 1.1 row of the version table
#      with "28.02.24"
#   3. Re-split a handful of runs ("Ameli" / "Pais" / "goods") at the exact
#      word boundaries that Word's background proofer later wraps in
#      <w:proofErr .../> (spell/grammar check flags). The COM object model
#      has no supported call that mints w:proofErr markers themselves (that
#      is purely an artifact of the interactive spell-checker), so this
#      script reproduces the one observable, automatable side effect of
#      that pass: the run boundaries it leaves behind, with formatting
#      carried through unchanged.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- helpers ---------------------------------------------------------------

# Returns a Range for the first occurrence of $needle at/after character
# offset $startAt, or $null if there is no further match.
function Find-OnceFrom($doc, $needle, $startAt) {
    $rng = $doc.Range($startAt, $doc.Content.End)
    $found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if (-not $found) { return $null }
    return $rng
}

# Returns an array of [start,end] pairs for every (non-overlapping)
# occurrence of $needle in the document.
function Find-All($doc, $needle) {
    $results = @()
    $searchStart = 0
    while ($true) {
        $rng = Find-OnceFrom $doc $needle $searchStart
        if ($null -eq $rng) { break }
        $results += , @($rng.Start, $rng.End)
        $searchStart = $rng.End
    }
    return $results
}

# Forces a run boundary immediately before absolute character offset $pos
# (so the text of the run is split in two there) without altering the
# visible formatting: flips Bold off/on (or on/off) across [$pos,$unitEnd)
# -- a real, recorded formatting operation -- then flips it straight back,
# which leaves Word no choice but to keep the run split while the net
# formatting is unchanged.
function Split-RunAt($doc, $pos, $unitEnd) {
    if ($pos -le 0 -or $pos -ge $unitEnd) { return }
    $probe = $doc.Range($pos, $unitEnd)
    $wasBold = $probe.Font.Bold
    if ($wasBold) {
        $probe.Font.Bold = $false
        $probe.Font.Bold = $true
    } else {
        $probe.Font.Bold = $true
        $probe.Font.Bold = $false
    }
}

# Splits the run(s) covering [$unitStart,$unitEnd) at every offset in
# $points (absolute character offsets). Processed right-to-left so each
# cut is made against an already-stable tail.
function Split-Unit($doc, $unitEnd, $points) {
    $sorted = $points | Sort-Object -Descending
    foreach ($p in $sorted) {
        Split-RunAt $doc $p $unitEnd
    }
}

# --- 1. "Group " + "J"  ->  "C#" -------------------------------------------

$d.Content.Find.Execute("Group J", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "C#", 2) | Out-Null

# --- 2. Version-table: blank Date cell for row "1.1" -> "28.02.24" --------

$table = $d.Tables.Item(1)
for ($r = 1; $r -le $table.Rows.Count; $r++) {
    $versionCell = $table.Cell($r, 1)
    if ($versionCell.Range.Text.TrimEnd([char]7, [char]13) -eq "1.1") {
        $dateCell = $table.Cell($r, 2)
        $dateCell.Range.Text = "28.02.24"
        break
    }
}

# --- 3. Re-split the runs the proofer later annotates ----------------------

# "Ameli Masewge Fernando" appears twice (the "Written by:" byline, and
# again in the Signature block); "Ameli Fernando" (bold, script font) is a
# separate single run right after the second one. In every case the split
# falls right after "Ameli" (5 characters in).
$masewges = Find-All $d "Ameli Masewge Fernando"
foreach ($p in $masewges) {
    $unitStart = $p[0]
    $unitEnd = $p[1]
    Split-Unit $d $unitEnd @($unitStart + "Ameli".Length)
}

$amelifernando = Find-OnceFrom $d "Ameli Fernando" 0
if ($null -ne $amelifernando) {
    $unitStart = $amelifernando.Start
    $unitEnd = $amelifernando.End
    Split-Unit $d $unitEnd @($unitStart + "Ameli".Length)
}

# "Nathan P. Pais " -> "Nathan P. " | "Pais" | " "
$nathan = Find-OnceFrom $d "Nathan P. Pais " 0
if ($null -ne $nathan) {
    $nStart = $nathan.Start
    $nEnd = $nathan.End
    $splitAfterPrefix = $nStart + ("Nathan P. ".Length)
    $splitAfterPais = $splitAfterPrefix + ("Pais".Length)
    Split-Unit $d $nEnd @($splitAfterPrefix, $splitAfterPais)
}

# ": Produced goods" -> ": Produced " | "goods"
$produced = Find-OnceFrom $d ": Produced goods" 0
if ($null -ne $produced) {
    $pStart = $produced.Start
    $pEnd = $produced.End
    $splitBeforeGoods = $pStart + (": Produced ".Length)
    Split-Unit $d $pEnd @($splitBeforeGoods)
}
